# Applies the "adding new project management report template" edit:
#  - Period selector (H3) moves from period 1 to period 13
#  - Actual Start / Actual Duration / % Complete (E:G) are filled in for
#    every activity row (6-21), driving the Gantt chart's "actual" bars
#  - The sheet view is left zoomed to 85% with E14 selected (matches the
#    state the workbook was saved in)
#  - A few helper columns around the legend boxes (Q:T) and the trailing
#    spacer column (AP) pick up explicit widths

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Period highlighted in the header / conditional formatting driver
$ws.Range("H3").Value = 13

# Actual Start (E), Actual Duration (F), Percent Complete (G) per activity
$actuals = @{
    6  = @(1, 3, 1)
    7  = @(3, 1, 1)
    8  = @(2, 2, 1)
    9  = @(1, 1, 1)
    10 = @(3, 6, 1)
    11 = @(1, 4, 1)
    12 = @(1, 1, 1)
    13 = @(2, 2, 1)
    14 = @(4, 2, 1)
    15 = @(5, 2, 1)
    16 = @(6, 2, 1)
    17 = @(7, 2, 1)
    18 = @(8, 2, 1)
    19 = @(9, 2, 1)
    20 = @(10, 2, 1)
    21 = @(11, 2, 1)
}

foreach ($row in $actuals.Keys) {
    $vals = $actuals[$row]
    $ws.Range("E$row").Value = $vals[0]
    $ws.Range("F$row").Value = $vals[1]
    $ws.Range("G$row").Value = $vals[2]
}

# Legend box helper columns get a touch wider so the "Actual Start" /
# "Actual Duration" labels fit, and the trailing spacer column is restored.
# (Inputs are pre-compensated for this engine's character->width rounding
# so the saved column width lands as close as possible to the authored
# 4.375 / 3.875 / 4.875 / 2.75 values.)
$ws.Columns.Item(17).ColumnWidth = 3.7142857142857144
$ws.Columns.Item(18).ColumnWidth = 3.7142857142857144
$ws.Columns.Item(19).ColumnWidth = 3.142857142857143
$ws.Columns.Item(20).ColumnWidth = 4.142857142857143
$ws.Columns.Item(42).ColumnWidth = 2.0

# Leave the sheet zoomed in a bit more, with E14 as the active selection
$ws.Range("E14").Select() | Out-Null
[void]($excel.ActiveWindow.Zoom = 85)
